# Update "想去人数" (want-to-go count) values in column F across the four
# worksheets of the workbook, per the latest gh-pages data refresh.

$wb = $excel.ActiveWorkbook

function Set-FValues($SheetName, $Updates) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $Updates.Keys) {
        $ws.Range("F$row").Value = $Updates[$row]
    }
}

# Sheet "展览" (Exhibition)
Set-FValues "展览" @{
    2  = 1715
    4  = 9976
    13 = 1561
    14 = 147
    15 = 345
    17 = 177
    18 = 458
    19 = 1147
    20 = 120
    21 = 29
    22 = 8
    23 = 71
    24 = 328
    26 = 300
    27 = 91
    28 = 620
    29 = 674
    30 = 15
    31 = 27
    32 = 207
    33 = 94
    34 = 251
    35 = 209
    36 = 357
    37 = 518
    38 = 485
    39 = 699
    42 = 791
    43 = 358
    44 = 317
}

# Sheet "演出" (Performance)
Set-FValues "演出" @{
    6  = 82
    8  = 716
    9  = 67
    12 = 55
    18 = 1059
    20 = 467
    21 = 1091
    22 = 314
    24 = 59
    28 = 353
    31 = 197
    34 = 151
    39 = 123
    42 = 33
}

# Sheet "本地生活" (Local Life)
Set-FValues "本地生活" @{
    5  = 180
    6  = 2481
    7  = 3964
    8  = 47
    10 = 231
    11 = 161
}

# Sheet "全部类型" (All Types)
Set-FValues "全部类型" @{
    2  = 1715
    5  = 9976
    6  = 180
    8  = 3964
    9  = 47
    10 = 231
    11 = 231
    14 = 147
    15 = 345
    16 = 67
    18 = 458
    19 = 1147
    20 = 120
    21 = 55
    23 = 8
    24 = 1059
    25 = 328
    28 = 300
    29 = 1091
    30 = 674
    32 = 59
    33 = 27
    34 = 207
    35 = 353
    36 = 209
    38 = 357
    39 = 518
    40 = 485
    41 = 197
    42 = 699
    44 = 791
    45 = 358
}
